$wb = $excel.ActiveWorkbook

# --- Sheet1 ("boaSignup"): add new test-case rows 11-15 in columns C/D ---
$ws1 = $wb.Worksheets.Item("boaSignup")

$ws1.Range("C11").Value = "FGHIJK"
$ws1.Range("D11").Value = "FHG1234$"

$ws1.Range("C12").Value = "LMNOPQ"
$ws1.Range("D12").Value = "pqrst123#"

$ws1.Range("C13").Value = "ABCDE2"
$ws1.Range("D13").Value = "abcd123#"

$ws1.Range("C14").Value = "FGHIJK2"
$ws1.Range("D14").Value = "FHG1234$"

$ws1.Range("C15").Value = "LMNO2"
$ws1.Range("D15").Value = "pqrst123#"

$ws1.Range("C11:D15").Select()

# --- Sheet2 ("loginNegativeTest"): remove ErrContains column and extra rows ---
$ws2 = $wb.Worksheets.Item("loginNegativeTest")

$ws2.Range("C1:C4").ClearContents()
$ws2.Range("A3:B4").ClearContents()

$ws2.Range("B1").Select()
